$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "58.792.92"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.518.80"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "535.80"
$ws.Range("E5").Value = "  +2.06%  "
Set-TextValue $ws.Range("D6") "136.20"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.565"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").Value = "2.519.58"
$ws.Range("E9").Value = "  +2.83%  "
Set-TextValue $ws.Range("D10") "0.101"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  -1.55%  "
Set-TextValue $ws.Range("D12") "5.36"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "2.965.06"
$ws.Range("E14").Value = "  +3.03%  "
Set-TextValue $ws.Range("D15") "23.00"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").Value = "58.780.06"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "2.526.73"
$ws.Range("E18").Value = "  +2.97%  "
Set-TextValue $ws.Range("D19") "11.07"
$ws.Range("E19").Value = "  +4.52%  "
Set-TextValue $ws.Range("D20") "4.26"
$ws.Range("E20").Value = "  +2.99%  "
Set-TextValue $ws.Range("D21") "322.78"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("E22").Value = "  +0.42%  "
Set-TextValue $ws.Range("D23") "5.96"
$ws.Range("E23").Value = "  +5.68%  "
Set-TextValue $ws.Range("D24") "65.09"
$ws.Range("E24").Value = "  +5.24%  "
Set-TextValue $ws.Range("D25") "0.420"
$ws.Range("E25").Value = "  +4.30%  "
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("E27").Value = "  +1.69%  "
Set-TextValue $ws.Range("D28") "7.51"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "0.0₃0768"
$ws.Range("E29").Value = "  +2.98%  "
Set-TextValue $ws.Range("D30") "6.64"
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.75"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D32") "170.27"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("E33").Value = "  +10.18%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +2.17%  "
Set-TextValue $ws.Range("D36") "18.34"
$ws.Range("E36").Value = "  +1.29%  "
Set-TextValue $ws.Range("D37") "4.05"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("E38").Value = "  +0.47%  "
Set-TextValue $ws.Range("D39") "36.83"
$ws.Range("E39").Value = "  +1.34%  "
Set-TextValue $ws.Range("D40") "0.808"
$ws.Range("E40").Value = "  +4.27%  "
$ws.Range("E41").Value = "  +1.99%  "
Set-TextValue $ws.Range("D42") "284.19"
$ws.Range("E42").Value = "  +4.98%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D43") "5.20"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "0.998"
$ws.Range("E44").Value = "  -0.24%  "
Set-TextValue $ws.Range("D45") "0.607"
$ws.Range("E45").Value = "  +4.18%  "
Set-TextValue $ws.Range("D46") "129.92"
$ws.Range("E46").Value = "  +8.86%  "
$ws.Range("E47").Value = "  +0.26%  "
Set-TextValue $ws.Range("D48") "0.0922"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  +1.00%  "
Set-TextValue $ws.Range("D51") "17.37"
$ws.Range("E51").Value = "  +3.14%  "

Write-Host "Applied 92 cell updates"
